$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra (blank) row 2, shifting the data row (row 3) up to row 2
# and the trailing blank-styled row (row 4) up to row 3.
$ws.Rows.Item(2).Delete()

# Update the active selection to match the post-edit state.
$ws.Range("D5").Select()
